# Insert a new product row ("TRIACTIN 4MG 20 TAB") into the pharmacy report,
# right after "SORAL 30 MG 30CAPS" (row 16) and before "VOLTAREN 75MG/3ML 3 AMP."
# (previously row 17). All following rows shift down by one; the running
# total and the summary/footer rows move down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 17 - this pushes rows 17.. down by one,
#    including the summary-total row and the footer row, and it carries the
#    formatting (styles, merges are NOT auto-copied so we redo them below).
$ws.Rows.Item(17).Insert()

# 2) Re-create the 3 merged ranges used by every product row, for the new row 17
$ws.Range("B17:G17").Merge()
$ws.Range("H17:K17").Merge()
$ws.Range("L17:M17").Merge()

# 3) Fill in the new row's data
$ws.Cells.Item(17, 1).Value = 14
$ws.Cells.Item(17, 2).Value = "TRIACTIN 4MG 20 TAB"
$ws.Cells.Item(17, 8).Value = "1:1"
$ws.Cells.Item(17, 12).Value = 23
$ws.Cells.Item(17, 14).Value = "0:2"

# 4) Fix row heights: the newly inserted row defaults to 14pt; restore it (and
#    its neighbours) to match the sheet's existing rhythm of alternating
#    24.75 / 25.5 pt rows.
$ws.Rows.Item(17).RowHeight = 25.5
$ws.Rows.Item(26).RowHeight = 24.75
$ws.Rows.Item(27).RowHeight = 26.25
$ws.Rows.Item(28).RowHeight = 16.5

# 5) Update the running-total cell (K27, formerly K26) to include the new
#    row's price (1533.04 + 23 = 1556.04).
$ws.Cells.Item(27, 11).Value = 1556.04
